$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: a new "white" entry under the searchSanPham/Dress column
$ws.Range("A3").Value = "white"

# Row 1: two new header cells for a simple login form
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "password"

# Row 6: another product row plus a sample login/contact record
$ws.Range("A6").Value = "Dress"
$ws.Range("D6").Value = 1234

# Email hyperlink cell in C6 (set the display text first so the hyperlink
# doesn't overwrite it with the raw mailto: address)
$ws.Range("C6").Value = "hanghang@gmai.com"
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:hanghang@gmai.com") | Out-Null

# Widen column C so the email address/headers are readable
$ws.Range("C1").ColumnWidth = 26.7

# Match the resulting selection
$ws.Range("A6").Select() | Out-Null
